$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab (was "Through 2021-09-06")
$ws.Name = "Through 2021-09-07"

# Update the "September (through 09-06)" label to "09-07"
$ws.Range("A10").Value = "September (through 09-07)"

# Update September row (row 10) values
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 17
$ws.Range("G10").Value = 24
$ws.Range("H10").Value = 28

# Update Total row (row 11) values
$ws.Range("B11").Value = 201
$ws.Range("C11").Value = 394
$ws.Range("D11").Value = 568
$ws.Range("E11").Value = 499
$ws.Range("F11").Value = 366
$ws.Range("G11").Value = 808
$ws.Range("H11").Value = 1099
